$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.158554560563064
$ws.Range("K2").Value = 0.0761309064869961
$ws.Range("L2").Value = 0.286786573395475
$ws.Range("N2").Value = 0.142637526134715
$ws.Range("B3").Value = 0.116851853607766
$ws.Range("K3").Value = 0.0725237685072292
$ws.Range("L3").Value = 0.251129621666901
$ws.Range("N3").Value = 0.104534936160411
$ws.Range("B4").Value = 0.0829083659428993
$ws.Range("K4").Value = 0.0197498970784084
$ws.Range("L4").Value = 0.0606228458463316
$ws.Range("N4").Value = 0.0942557965998385
$ws.Range("B5").Value = 0.0714895172677557
$ws.Range("F5").Value = 0.0624489863806241
$ws.Range("K5").Value = 0.040688726015579
$ws.Range("L5").Value = 0.0431192561344624
$ws.Range("N5").Value = 0.0842750030357131
$ws.Range("B6").Value = 0.05450527325568
$ws.Range("D6").Value = 0.0718687212048466
$ws.Range("K6").Value = 0.0129576420740687
$ws.Range("L6").Value = 0.0259958246698663
$ws.Range("N6").Value = 0.0627440974232051
$ws.Range("B7").Value = 0.0455409606524207
$ws.Range("I7").Value = 0.0487964302518536
$ws.Range("K7").Value = 0.0160044958444441
$ws.Range("L7").Value = 0.0347523618529367
$ws.Range("N7").Value = 0.0498049899179879
$ws.Range("B8").Value = 0.0387328434737347
$ws.Range("K8").Value = 0.0105116378040686
$ws.Range("L8").Value = 0.0282733934651745
$ws.Range("N8").Value = 0.0456258132835108
